$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 190 — all existing rows from 190 downward
# (190-239) shift down to (191-240), matching the target workbook's extra
# weekly price-report entry for Pomelo "Start Ruby" / "Primera".
$ws.Rows.Item(190).Insert()

$ws.Range("A190").Value = 10
$ws.Range("B190").Value = "Vega Modelo de Temuco"
$ws.Range("C190").Value = "La Araucanía"
$ws.Range("D190").Value = 44722
$ws.Range("E190").Value = 9
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100102
$ws.Range("H190").Value = "Cítricos"
$ws.Range("I190").Value = 100102006
$ws.Range("J190").Value = "Pomelo"
$ws.Range("K190").Value = "Start Ruby"
$ws.Range("L190").Value = "Primera"
$ws.Range("M190").Value = 55
$ws.Range("N190").Value = 12000
$ws.Range("O190").Value = 12000
$ws.Range("P190").Value = 12000
$ws.Range("Q190").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R190").Value = "Región de O'Higgins"
$ws.Range("S190").Value = 800
$ws.Range("T190").Value = 15
